$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SUPPORT API")

$ws.Range("J2").Value = "status=200||access_token=00D3000000000un!AR8AQC7cJ81UAeuz.9pdSgmblAlTzFyJzZpIiIOvRwg..r3e.y4142nomU4aNdk.xoIopCgcmPiZvQy8EHfsOPVwykB1_AKq||instance_url=https://na33.salesforce.com"
